$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in BOM: J4 (OLED module) row - MOUSER and DIGIKEY columns were
# showing "-" (not applicable) instead of "N.M." (not mounted), which is
# inconsistent with the MPN/LCSC columns on the same row and with the
# MK1/MK2 microphone row below it.
$ws.Range("H9").Value = "N.M."
$ws.Range("I9").Value = "N.M."

# Update the last active selection left over from editing (reflects the
# modules-pictures work done elsewhere in the sheet before saving).
$ws.Range("L16").Select()
